$d = $word.ActiveDocument

function Merge-ParagraphRuns {
    param(
        [int]$ParaIndex,
        [int]$KeepLen,
        [string]$FullText
    )

    $p = $d.Paragraphs($ParaIndex)
    $paraRange = $p.Range
    $paraStart = $paraRange.Start
    $paraEnd = $paraRange.End - 1   # exclude the trailing paragraph mark

    # Keep the paragraph's first run text node untouched (so its existing
    # xml:space="preserve" attribute survives) and delete everything that
    # follows it, then append the remaining text with InsertAfter so the
    # new text is merged into that same run instead of creating new runs.
    $keepEnd = $paraStart + $KeepLen
    $restRange = $d.Range($keepEnd, $paraEnd)
    $restRange.Delete()

    $insPoint = $d.Range($keepEnd, $keepEnd)
    $remainder = $FullText.Substring($KeepLen)
    $insPoint.InsertAfter($remainder)
}

# Title: "Questions:" + " Trigonometry (radians)"
Merge-ParagraphRuns 1 10 "Questions: Trigonometry (radians)"

# Author: "Dzhemma" + " Ruseva, Ellie Gurini, Ciara Cormican"
Merge-ParagraphRuns 2 7 "Dzhemma Ruseva, Ellie Gurini, Ciara Cormican"

# Abstract: "A" + " selection of questions on trigonometry, where angles are measured in degrees."
Merge-ParagraphRuns 4 1 "A selection of questions on trigonometry, where angles are measured in degrees."
